# Rename the worksheet to reflect the user story it now contains.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Foglio 1 - User Story"

# Replace the placeholder user-story text in A2 with the real one.
$ws.Cells.Item(2, 1).Value = "As a radiologist, I want to use the ID3 algorithm to develop decision tree models for diagnosing and predicting medical conditions based on various medical imaging data, such as X-rays, CT scans, and MRI scans.`n"

# Widen column A so the longer text fits better, and grow the header rows.
$ws.Columns.Item(1).ColumnWidth = 119 - 5/7
$ws.Rows.Item(2).RowHeight = 44.2
$ws.Rows.Item(3).RowHeight = 20.25

# Freeze the first two (header) rows.
$ws.Range("A3").Select()
$excel.ActiveWindow.FreezePanes = $true

# Extend the table with one more blank row, matching the formatting of the
# existing blank rows below the header.
$ws.Cells.Item(23, 1).NumberFormat = "General"
